$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = "sfhsdljfds"

$ws.Range("E12").Select()
